$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell F1 - matches style of other header cells (B1:E1)
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "time_taken"

$timestamps = @(
    "2021-10-05 10:50:39.302590",
    "2021-10-05 10:50:39.302600",
    "2021-10-05 10:50:39.302603",
    "2021-10-05 10:50:39.302606",
    "2021-10-05 10:50:39.302609",
    "2021-10-05 10:50:39.302611",
    "2021-10-05 10:50:39.302614",
    "2021-10-05 10:50:39.302616",
    "2021-10-05 10:50:39.302619",
    "2021-10-05 10:50:39.302622",
    "2021-10-05 10:50:39.302624",
    "2021-10-05 10:50:39.302627",
    "2021-10-05 10:50:39.302629",
    "2021-10-05 10:50:39.302631",
    "2021-10-05 10:50:39.302634",
    "2021-10-05 10:50:39.302636",
    "2021-10-05 10:50:39.302639",
    "2021-10-05 10:50:39.302642",
    "2021-10-05 10:50:39.302644",
    "2021-10-05 10:50:39.302647",
    "2021-10-05 10:50:39.302649",
    "2021-10-05 10:50:39.302652",
    "2021-10-05 10:50:39.302654",
    "2021-10-05 10:50:39.302657",
    "2021-10-05 10:50:39.302659",
    "2021-10-05 10:50:39.302662",
    "2021-10-05 10:50:39.302665",
    "2021-10-05 10:50:39.302667"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
